# switched to yahoo finance
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantities for existing holdings (prices refreshed from Yahoo Finance)
$ws.Range("C2").Value = 300
$ws.Range("C5").Value = 200

# Add a new holding row: TTD (formatting matches the row above)
$ws.Range("A9:F9").Copy()
$ws.Range("A10").PasteSpecial(-4122)

$ws.Range("A10").Value = 43503
$ws.Range("B10").Value = "TTD"
$ws.Range("C10").Value = 100
$ws.Range("D10").Value = 14.3
$ws.Range("E10").Formula = "=+C10*D10"
$ws.Range("F10").Value = 44196

$ws.Range("C3").Select()

$wb.Save()
